$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Galaxy S5 SM-G900A ShirNate" -> add googleTest row, result FAIL, with screenshot hyperlink
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A5").Value = "googleTest"
$ws1.Range("C5").Value = "FAIL"
$ws1.Range("C5").Interior.ColorIndex = 3
$ws1.Hyperlinks.Add($ws1.Range("C5"), "C%3A%2FUsers%2FAvnerG%2Fgit%2FBeton%2FBeton%2Ftest-output%2Fscreenshots-tests%2F2015-08-26-11-37-00-IDT.png")

# --- Sheet 2: "iPhone-6 Avner" -> add googleTest row, result PASS
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A5").Value = "googleTest"
$ws2.Range("C5").Value = "PASS"
$ws2.Range("C5").Interior.ColorIndex = 10

# --- Sheet 3: "iPhone-6 Raj" -> add googleTest row, result PASS
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A5").Value = "googleTest"
$ws3.Range("C5").Value = "PASS"
$ws3.Range("C5").Interior.ColorIndex = 10
